# Linking settings tab to selection tab
#
# - Rename the three sheets to their properly-cased names.
# - Move the "active"/selected tab from the English sheet to the Exiobase
#   (settings) sheet, and update each sheet's remembered cell selection.

$wb = $excel.ActiveWorkbook

$wsExiobase = $wb.Worksheets.Item(1)
$wsGerman   = $wb.Worksheets.Item(2)
$wsEnglish  = $wb.Worksheets.Item(3)

# Rename sheets (exiobase -> Exiobase, german -> Deutsch, english -> English)
$wsExiobase.Name = "Exiobase"
$wsGerman.Name   = "Deutsch"
$wsEnglish.Name  = "English"

# Update stored selections on the non-active sheets first...
$wsGerman.Range("F27").Select()
$wsEnglish.Range("F34").Select()

# ...then activate the Exiobase sheet and set its selection, so it becomes
# the workbook's active/selected tab (moving tabSelected off of English).
$wsExiobase.Activate()
$wsExiobase.Range("C34").Select()
